$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade record in row 4 (mirrors the existing row 3 trade entry)
$ws.Range("A4").Value = 42633.676701388889
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9956.5
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 108.67
$ws.Range("F4").Value = 107.73
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -0.87
$ws.Range("I4").Value = $false

# Row 3's date column (A) carries a date/time number format; copy that
# formatting onto the new date cell so it reuses the existing style (rather
# than Excel minting a brand new one just for this cell).
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
